$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1932.6666
$ws.Range("I17").Value = 1800
$ws.Range("K17").Value = 5400
$ws.Range("M17").Value = -5232
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H52").Value = 4250
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 4250
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 12750
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -13070
$ws.Range("H64").Value = 2964.2144
$ws.Range("I64").Value = 2916.6667
$ws.Range("J64").Value = 3249.5
$ws.Range("K64").Value = 2916.6667
$ws.Range("L64").Value = 3249.5
$ws.Range("M64").Value = -2668.6667
$ws.Range("N64").Value = -3745.5
$ws.Range("H67").Value = 2964.2144
$ws.Range("I67").Value = 2916.6667
$ws.Range("J67").Value = 3249.5
$ws.Range("K67").Value = 2916.6667
$ws.Range("L67").Value = 3249.5
$ws.Range("M67").Value = -2058.6667
$ws.Range("N67").Value = -4965.5
$ws.Range("H80").Value = 433.375
$ws.Range("I80").Value = 448.33334
$ws.Range("J80").Value = 424.4
$ws.Range("K80").Value = 1345.00002
$ws.Range("L80").Value = 1273.2
$ws.Range("M80").Value = -347.0000199999999
$ws.Range("N80").Value = -3269.2
$ws.Range("H83").Value = 433.375
$ws.Range("I83").Value = 448.33334
$ws.Range("J83").Value = 424.4
$ws.Range("K83").Value = 4035.00006
$ws.Range("L83").Value = 3819.6
$ws.Range("M83").Value = 956.9999399999997
$ws.Range("N83").Value = -13803.6
$ws.Range("H103").Value = 740
$ws.Range("I103").Value = 300
$ws.Range("J103").Value = 850
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 2550
$ws.Range("N103").Value = -3722
$ws.Range("M103").Value = -314
$ws.Range("H111").Value = 441
$ws.Range("I111").Value = 232.33333
$ws.Range("J111").Value = 649.6667
$ws.Range("K111").Value = 696.99999
$ws.Range("L111").Value = 1949.0001
$ws.Range("M111").Value = 2370.00001
$ws.Range("N111").Value = -8083.0001
$ws.Range("H125").Value = 2999
$ws.Range("I125").Value = 2999
$ws.Range("K125").Value = 26991
$ws.Range("M125").Value = -24531
$ws.Range("H135").Value = 1255.75
$ws.Range("I135").Value = 1363.8572
$ws.Range("K135").Value = 12274.7148
$ws.Range("M135").Value = -9739.7148
$ws.Range("H138").Value = 3333.1785
$ws.Range("I138").Value = 1407.1428
$ws.Range("J138").Value = 3975.1904
$ws.Range("K138").Value = 4221.428400000001
$ws.Range("L138").Value = 11925.5712
$ws.Range("M138").Value = 918.5715999999993
$ws.Range("N138").Value = -22205.5712
$ws.Range("H141").Value = 35655.715
$ws.Range("I141").Value = 35655.715
$ws.Range("K141").Value = 106967.145
$ws.Range("M141").Value = -101787.145
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H45").Value = 2365
$ws.Range("I45").Value = 2389.4
$ws.Range("K45").Value = 2389.4
$ws.Range("M45").Value = -2012.4
$ws.Range("H63").Value = 3500
$ws.Range("I63").Value = 3500
$ws.Range("K63").Value = 3500
$ws.Range("M63").Value = -2814
$ws.Range("H66").Value = 3500
$ws.Range("I66").Value = 3500
$ws.Range("K66").Value = 17500
$ws.Range("M66").Value = -14068
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H132").Value = 1573
$ws.Range("I132").Value = 1573
$ws.Range("K132").Value = 4719
$ws.Range("M132").Value = -2189
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H134").Value = 763
$ws.Range("I134").Value = 763
$ws.Range("K134").Value = 2289
$ws.Range("M134").Value = 246
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10775.667
$ws.Range("I31").Value = 10775.667
$ws.Range("K31").Value = 10775.667
$ws.Range("M31").Value = -10480.667
$ws.Range("H34").Value = 10775.667
$ws.Range("I34").Value = 10775.667
$ws.Range("K34").Value = 10775.667
$ws.Range("M34").Value = -10573.667
$ws.Range("H58").Value = 2262.6667
$ws.Range("I58").Value = 2194.8572
$ws.Range("K58").Value = 2194.8572
$ws.Range("M58").Value = -1991.8572
$ws.Range("H62").Value = 11288.1
$ws.Range("I62").Value = 14998.8
$ws.Range("K62").Value = 14998.8
$ws.Range("M62").Value = -14374.8
$ws.Range("H65").Value = 11288.1
$ws.Range("I65").Value = 14998.8
$ws.Range("K65").Value = 74994
$ws.Range("M65").Value = -71874
$ws.Range("H136").Value = 2262.6667
$ws.Range("I136").Value = 2194.8572
$ws.Range("K136").Value = 6584.571599999999
$ws.Range("M136").Value = -4034.571599999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1057.875
$ws.Range("J92").Value = 395
$ws.Range("L92").Value = 1185
$ws.Range("N92").Value = -3681
$ws.Range("H122").Value = 4354.75
$ws.Range("J122").Value = 4354.75
$ws.Range("L122").Value = 39192.75
$ws.Range("N122").Value = -44092.75
$ws.Range("H129").Value = 8166.6665
$ws.Range("I129").Value = 8166.6665
$ws.Range("K129").Value = 24499.9995
$ws.Range("M129").Value = -19499.9995
$ws.Range("H132").Value = 2667.3333
$ws.Range("I132").Value = 2751
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 24759
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -22229
$ws.Range("N132").Value = -27560
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9011.125
$ws.Range("I81").Value = 9727.143
$ws.Range("J81").Value = 3999
$ws.Range("K81").Value = 19454.286
$ws.Range("L81").Value = 7998
$ws.Range("M81").Value = -18393.286
$ws.Range("N81").Value = -10120
$ws.Range("H84").Value = 9011.125
$ws.Range("I84").Value = 9727.143
$ws.Range("J84").Value = 3999
$ws.Range("K84").Value = 97271.42999999999
$ws.Range("L84").Value = 39990
$ws.Range("M84").Value = -91967.42999999999
$ws.Range("N84").Value = -50598
$ws.Range("H136").Value = 2726.7273
$ws.Range("I136").Value = 2726.7273
$ws.Range("K136").Value = 8180.1819
$ws.Range("M136").Value = -5630.1819
